# Apply tissue-DB fixes to sln_tissue__FluorescentReagent sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize casing on a couple of reagent_name values ---
$ws.Range("A2").Value2 = "streptavidin_568"
$ws.Range("A3").Value2 = "streptavidin_488"

# --- Disambiguate donkey secondary antibodies by vendor ---
$ws.Range("A9").Value2  = "donkey-X-goat_647_jackson"
$ws.Range("A12").Value2 = "donkey-X-mouse_488_jackson"
$ws.Range("A22").Value2 = "donkey-X-mouse_488_thermofisher"
$ws.Range("A31").Value2 = "donkey-X-goat_647_thermofisher"

# --- Add a new DAPI reagent row ---
$ws.Range("A32").Value2 = "DAPI"
$ws.Range("B32").Value2 = "4′,6-Diamidine-2′-phenylindole dihydrochloride"
$ws.Range("C32").Value2 = 10236276001
$ws.Range("D32").Value2 = 2000
$ws.Range("E32").Value2 = "NA"
$ws.Range("F32").Value2 = "NA"
$ws.Range("G32").Value2 = "NA"
$ws.Range("H32").Value2 = "DAPI"
$ws.Range("I32").Value2 = "Millipore Sigma"

# --- Copy number/border formatting onto the new row from existing rows ---
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null

$ws.Range("I2").Copy() | Out-Null
$ws.Range("I32").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Match the saved selection state from the diff ---
$ws.Range("C17").Select()
